# Insert a new data row at row 680 (pushing the existing rows 680-759 down
# to 681-760) and populate the new row with the latest price-survey entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("680:680").Insert()

$ws.Range("A680").Value = 4
$ws.Range("B680").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C680").Value = "Los Lagos"
$ws.Range("D680").Value = 44918
$ws.Range("E680").Value = 10
$ws.Range("F680").Value = 100112004
$ws.Range("G680").Value = "Cebolla"
$ws.Range("H680").Value = "Sin especificar"
$ws.Range("I680").Value = "1a nueva(o)"
$ws.Range("J680").Value = 700
$ws.Range("K680").Value = 16000
$ws.Range("L680").Value = 17000
$ws.Range("M680").Value = 16500
$ws.Range("N680").Value = "$/malla 18 kilos"
$ws.Range("O680").Value = "Región de O'Higgins"
$ws.Range("P680").Value = 917
$ws.Range("Q680").Value = 18
$ws.Range("R680").Value = "Hortaliza"
